# Apply the weekly cryptocurrency data refresh to Sheet1 (cryptos.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.398.18"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.880.78"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.46"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07942"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3139"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "1.892.59"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.77"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7073"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.394"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008409"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "29.437.36"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.69"
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.34"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "2.142.89"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.675"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.052"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.62"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.92"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.311"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.235"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05301"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.918"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7552"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.705"
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.283.46"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01879"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.767"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.369"
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9062"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.20"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000131"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "2.040.82"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.804"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5211"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.503"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4340"
$ws.Range("E51").Value = "  -0.57%  "
